$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.874.79"
$ws.Range("E2").Value = "  -2.28%  "

# Row 3
$ws.Range("D3").Value = "1.754.48"
$ws.Range("E3").Value = "  -4.64%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.48"
$ws.Range("E5").Value = "  -8.12%  "

# Row 6
$ws.Range("E6").Value = "  +0.01%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5088"
$ws.Range("E7").Value = "  -5.23%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.42"
$ws.Range("E8").Value = "  -5.31%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2766"
$ws.Range("E9").Value = "  -6.20%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06198"
$ws.Range("E10").Value = "  -10.82%  "

# Row 11
$ws.Range("D11").Value = "1.752.50"
$ws.Range("E11").Value = "  -5.41%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06970"
$ws.Range("E12").Value = "  -3.06%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.71"
$ws.Range("E13").Value = "  -9.20%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6136"
$ws.Range("E14").Value = "  -15.59%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.531"
$ws.Range("E15").Value = "  -8.96%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.45"
$ws.Range("E16").Value = "  -13.13%  "

# Row 17
$ws.Range("E17").Value = "  +0.03%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  -0.01%  "

# Row 19
$ws.Range("D19").Value = "25.878.06"
$ws.Range("E19").Value = "  -2.32%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006885"
$ws.Range("E20").Value = "  -12.79%  "

# Row 21
$ws.Range("E21").Value = "  -15.04%  "

# Row 22
$ws.Range("D22").Value = "1.972.05"
$ws.Range("E22").Value = "  -5.28%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.081"
$ws.Range("E23").Value = "  -11.01%  "

# Row 24
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.257"
$ws.Range("E24").Value = "  -9.93%  "

# Row 25
$ws.Range("B25").Value = "Chainlink"
$ws.Range("C25").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.251"
$ws.Range("E25").Value = "  -12.24%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "137.97"
$ws.Range("E26").Value = "  -3.48%  "

# Row 27
$ws.Range("E27").Value = "  -12.70%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.08"
$ws.Range("E28").Value = "  -11.07%  "

# Row 29
$ws.Range("E29").Value = "  -15.66%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "103.71"
$ws.Range("E30").Value = "  -6.48%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08251"
$ws.Range("E31").Value = "  -7.30%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.704"
$ws.Range("E32").Value = "  -13.00%  "

# Row 33
$ws.Range("E33").Value = "  -13.26%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04570"
$ws.Range("E34").Value = "  -5.48%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9999"
$ws.Range("E35").Value = "  +0.00%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.641"
$ws.Range("E36").Value = "  -9.05%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9940"
$ws.Range("E37").Value = "  -12.11%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6113"
$ws.Range("E38").Value = "  -15.72%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.704"
$ws.Range("E39").Value = "  -12.63%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01560"
$ws.Range("E40").Value = "  -8.65%  "

# Row 41
$ws.Range("E41").Value = "  +0.04%  "

# Row 42
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.901"
$ws.Range("E42").Value = "  -16.85%  "

# Row 43
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.56"
$ws.Range("E43").Value = "  -3.54%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3856"
$ws.Range("E44").Value = "  -17.31%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7414"
$ws.Range("E45").Value = "  -17.77%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.967"
$ws.Range("E46").Value = "  -15.45%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05435"
$ws.Range("E47").Value = "  -5.26%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1114"
$ws.Range("E48").Value = "  -10.33%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.026"
$ws.Range("E49").Value = "  -18.68%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.703"
$ws.Range("E50").Value = "  -14.40%  "

# Row 51
$ws.Range("E51").Value = "  -13.56%  "
